$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text ("@") number format on the cells being updated so that
# numeric- and percent-looking strings are stored as literal text,
# matching the inline string values from the source data feed.
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "B17", "C17", "D17", "E17", "B18", "C18", "D18", "E18", "B19", "C19", "D19", "E19", "B20", "C20", "D20", "E20", "B21", "C21", "D21", "E21", "B22", "C22", "D22", "E22", "B23", "C23", "D23", "E23", "B24", "C24", "D24", "E24", "D25", "E25", "D26", "E26", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "296.81"
$ws.Range("E2").Value = "2.38%"
$ws.Range("D3").Value = "40.62"
$ws.Range("E3").Value = "3.04%"
$ws.Range("D4").Value = "5.055"
$ws.Range("E4").Value = "0.49%"
$ws.Range("D5").Value = "0.07435"
$ws.Range("E5").Value = "1.34%"
$ws.Range("D6").Value = "4.356"
$ws.Range("E6").Value = "1.69%"
$ws.Range("D7").Value = "1.572"
$ws.Range("E7").Value = "1.52%"
$ws.Range("D8").Value = "0.9352"
$ws.Range("E8").Value = "2.55%"
$ws.Range("D9").Value = "2.402"
$ws.Range("E9").Value = "0.21%"
$ws.Range("E10").Value = "1.04%"
$ws.Range("D11").Value = "0.1807"
$ws.Range("E11").Value = "3.76%"
$ws.Range("D12").Value = "0.08819"
$ws.Range("E12").Value = "1.39%"
$ws.Range("D13").Value = "0.04305"
$ws.Range("E13").Value = "3.55%"
$ws.Range("D14").Value = "0.1047"
$ws.Range("E14").Value = "-0.52%"
$ws.Range("D15").Value = "0.001263"
$ws.Range("E15").Value = "-0.87%"
$ws.Range("D16").Value = "0.005956"
$ws.Range("E16").Value = "2.40%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.357"
$ws.Range("E17").Value = "-1.15%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "0.3307"
$ws.Range("E18").Value = "0.68%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "7.929"
$ws.Range("E19").Value = "4.76%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1376"
$ws.Range("E20").Value = "1.86%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "0.2960"
$ws.Range("E21").Value = "2.68%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "0.04015"
$ws.Range("E22").Value = "4.63%"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "0.001267"
$ws.Range("E23").Value = "-0.09%"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "0.003868"
$ws.Range("E24").Value = "-0.62%"
$ws.Range("D25").Value = "0.0001227"
$ws.Range("E25").Value = "-4.30%"
$ws.Range("D26").Value = "0.0003713"
$ws.Range("E26").Value = "-0.39%"
$ws.Range("D38").Value = "0.02396"
$ws.Range("E38").Value = "2.64%"
$ws.Range("D39").Value = "0.05197"
$ws.Range("E39").Value = "3.48%"
$ws.Range("D40").Value = "0.005954"
$ws.Range("E40").Value = "16.53%"
$ws.Range("D41").Value = "0.007765"
$ws.Range("E41").Value = "0.85%"
$ws.Range("D42").Value = "0.1320"
$ws.Range("E42").Value = "3.73%"
$ws.Range("D43").Value = "0.007364"
$ws.Range("E43").Value = "-0.11%"
$ws.Range("D44").Value = "0.007188"
$ws.Range("E44").Value = "3.12%"
$ws.Range("D45").Value = "0.2952"
$ws.Range("E45").Value = "-6.10%"
$ws.Range("D46").Value = "0.00006263"
$ws.Range("E46").Value = "-3.87%"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").Value = "-0.38%"
$ws.Range("D48").Value = "0.04552"
$ws.Range("E48").Value = "-81.92%"
$ws.Range("D49").Value = "0.004190"
$ws.Range("E49").Value = "-0.38%"
$ws.Range("D50").Value = "0.00002095"
$ws.Range("E50").Value = "-0.38%"
$ws.Range("D51").Value = "0.0001995"
$ws.Range("E51").Value = "-0.38%"
